# Add a new "2022-Q1" sheet (positioned after "2021-Q4", before "总计")
# and prepend a matching summary row on the "总计" sheet.
# See commit message: "feat: add 2022-Q1 data"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q1" worksheet right after "2021-Q4"
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add($null, $q4)
$newSheet.Name = "2022-Q1"

# Reuse the exact header formatting (bold / centered / boxed) and the
# index-column formatting from the "2021-Q4" sheet, which uses the same
# column layout.
$q4.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$q4.Range("A2").Copy()
$newSheet.Range("A2:A4").PasteSpecial(-4122)

# Header row text
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $newSheet.Cells.Item(1, $i + 2).Value = $headers[$i]
}

# Fund holdings data rows: index, code, name, size, stock position,
# position ratio, held value (100M CNY), position rank
$rows = @(
    @(0, "160642", "鹏华增瑞灵活配置混合(LOF)", "6.76", "91.34", "9.44", "0.6381", 2),
    @(1, "001675", "江信同福灵活配置混合A",     "1.02", "92.43", "6.28", "0.0641", 8),
    @(2, "001676", "江信同福灵活配置混合C",     "0.42", "92.43", "6.28", "0.0264", 8)
)

foreach ($row in $rows) {
    $r = [int]$row[0] + 2

    $newSheet.Cells.Item($r, 1).Value = $row[0]

    $bCell = $newSheet.Cells.Item($r, 2)
    $bCell.NumberFormat = "@"
    $bCell.Value = $row[1]

    $newSheet.Cells.Item($r, 3).Value = $row[2]

    $dCell = $newSheet.Cells.Item($r, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $row[3]

    $eCell = $newSheet.Cells.Item($r, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value = $row[4]

    $fCell = $newSheet.Cells.Item($r, 6)
    $fCell.NumberFormat = "@"
    $fCell.Value = $row[5]

    $gCell = $newSheet.Cells.Item($r, 7)
    $gCell.NumberFormat = "@"
    $gCell.Value = $row[6]

    $newSheet.Cells.Item($r, 8).Value = $row[7]
}

# ---------------------------------------------------------------------
# 2) Update the "总计" sheet: insert a new summary row for 2022-Q1
#    above the existing 2021-Q4 / 2021-Q3 rows.
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

# Keep a formatting reference to the (soon to be pushed down) first
# data row before we insert, so the new row 2 can reuse its index-column
# style exactly.
$wsTotal.Rows.Item(2).Insert()
$wsTotal.Range("A2:D2").ClearFormats()

$wsTotal.Cells.Item(3, 1).Copy()
$wsTotal.Cells.Item(2, 1).PasteSpecial(-4122)

$wsTotal.Cells.Item(2, 1).Value = 0
$wsTotal.Cells.Item(2, 2).Value = "2022-Q1"
$wsTotal.Cells.Item(2, 3).Value = 3
$wsTotal.Cells.Item(2, 4).Value = 0.73

# Renumber the 0-based index column for the rows pushed down
$wsTotal.Cells.Item(3, 1).Value = 1
$wsTotal.Cells.Item(4, 1).Value = 2

Write-Host "2022-Q1 sheet added and summary sheet updated"
